# Commit: update hotel reviews data
# - hotel_info (sheet1) row 2: fill English_Reviews_num / Local_Rank /
#   Total_Reviews_num / Orbitz_ReviewURL, which were blank placeholders.
# - review_info (sheet2): add 3 review rows (rows 2-4).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("hotel_info")
$ws2 = $wb.Worksheets.Item("review_info")

# Write a cell as TEXT even when the content looks numeric/date-like.
# A plain .Value assignment mirrors typing into Excel, so "7", "08/11/2018",
# "July 2014", etc. would silently become a number/date. Prefixing the
# literal with a quote forces text entry (matching the source data, which
# stores these as strings); resetting .Style back to "Normal" afterwards
# clears the quote-prefix formatting again so no stray style is left on
# the cell.
function Set-TextCell {
    param($ws, $row, $col, $text)
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

function Set-NumCell {
    param($ws, $row, $col, $num)
    $ws.Cells.Item($row, $col).Value = $num
}

# --- hotel_info: row 2, columns G (English_Reviews_num), H (Local_Rank),
#     I (Total_Reviews_num), J (Orbitz_ReviewURL) ---
Set-TextCell $ws1 2 7 '7'
Set-TextCell $ws1 2 8 '15'
Set-TextCell $ws1 2 9 '7'
Set-TextCell $ws1 2 10 '?'

# --- review_info: row 2 ---
Set-NumCell $ws2 2 1 34156
Set-NumCell $ws2 2 4 1
Set-TextCell $ws2 2 5 '08/11/2018'
Set-TextCell $ws2 2 6 'https://www.tripadvisor.com/ShowUserReviews-g56855-d223186-r238551340-InTown_Suites_NASA-Webster_Texas.html'
Set-TextCell $ws2 2 7 '56855'
Set-TextCell $ws2 2 8 '223186'
Set-TextCell $ws2 2 9 '238551340'
Set-TextCell $ws2 2 10 '11/06/2014'
Set-TextCell $ws2 2 11 'Not clean'
Set-TextCell $ws2 2 12 'I reserved this hotel to stay in while I was on assignment.  I am glad I looked at the room before I rented.  There was a great deal of mold in the bathroom, the carpet was disgustingly dirty and the smell was horrible. If I could have rated it a zero I would have. '
Set-NumCell $ws2 2 13 1
Set-NumCell $ws2 2 22 0
Set-TextCell $ws2 2 25 'I reserved this hotel to stay in while I was on assignment.  I am glad I looked at the room before I rented.  There was a great deal of mold in the bathroom, the carpet was disgustingly dirty and the smell was horrible. If I could have rated it a zero I would have. '

# --- review_info: row 3 ---
Set-NumCell $ws2 3 1 34156
Set-NumCell $ws2 3 4 2
Set-TextCell $ws2 3 5 '08/11/2018'
Set-TextCell $ws2 3 6 'https://www.tripadvisor.com/ShowUserReviews-g56855-d223186-r213559603-InTown_Suites_NASA-Webster_Texas.html'
Set-TextCell $ws2 3 7 '56855'
Set-TextCell $ws2 3 8 '223186'
Set-TextCell $ws2 3 9 '213559603'
Set-TextCell $ws2 3 10 '07/03/2014'
Set-TextCell $ws2 3 11 'she forced her way into our room'
Set-TextCell $ws2 3 12 'I wish there was a negative star rating. The manager, Robyn Nelson (---)799-9527, forced her way 3 times into our room while my husband and I were naked. She used extreme profanity while throwing our belongings around, flipped the bed in the air and screaming the entire time. Previously she approached me about joining a pyramid scam selling coffee and soap. I was non committal about attending a 4th of July party given by her ''mentor'' friend, I was polite but could barely get my laundry card without her telling me she needed 8 more people so she could get her Mercedes. She was violent, abusive and destructive. She believed we had a cat in the room and not only went though our belongings when we were not there, she showed up at 10:30 am demanding to search us. We asked several other locations to advise us, and they handled it professionally, advising us to call the corporate office. My husband spoke to Mr. Carter, who documented everything and told us he had 3 other calls to make regarding this individual. He told him that most likely, she would not be there when we got back from work. He was polite and apologetic. He reassured us that we would not have to be afraid to return to our room. He was also very eager to make us feel safe, welcomed, and comfortable at Inn Town Suites.MoreShow less'
Set-NumCell $ws2 3 13 1
Set-TextCell $ws2 3 14 'July 2014'
Set-TextCell $ws2 3 15 ' traveled on business'
Set-NumCell $ws2 3 22 0
Set-TextCell $ws2 3 25 'I wish there was a negative star rating. The manager, Robyn Nelson (---)799-9527, forced her way 3 times into our room while my husband and I were naked. She used extreme profanity while throwing our belongings around, flipped the bed in the air and screaming the entire time. Previously she approached me about joining a pyramid scam selling coffee and soap. I was non committal about attending a 4th of July party given by her ''mentor'' friend, I was polite but could barely get my laundry card without her telling me she needed 8 more people so she could get her Mercedes. She was violent, abusive and destructive. She believed we had a cat in the room and not only went though our belongings when we were not there, she showed up at 10:30 am demanding to search us. We asked several other locations to advise us, and they handled it professionally, advising us to call the corporate office. My husband spoke to Mr. Carter, who documented everything and told us he had 3 other calls to make regarding this individual. He told him that most likely, she would not be there when we got back from work. He was polite and apologetic. He reassured us that we would not have to be afraid to return to our room. He was also very eager to make us feel safe, welcomed, and comfortable at Inn Town Suites.More'

# --- review_info: row 4 ---
Set-NumCell $ws2 4 1 34156
Set-NumCell $ws2 4 4 3
Set-TextCell $ws2 4 5 '08/11/2018'
Set-TextCell $ws2 4 6 'https://www.tripadvisor.com/ShowUserReviews-g56855-d223186-r186225775-InTown_Suites_NASA-Webster_Texas.html'
Set-TextCell $ws2 4 7 '56855'
Set-TextCell $ws2 4 8 '223186'
Set-TextCell $ws2 4 9 '186225775'
Set-TextCell $ws2 4 10 '11/29/2013'
Set-TextCell $ws2 4 11 'Fairly Nice But...'
Set-TextCell $ws2 4 12 'I lived here 6 months. It is a great location---the management is nice. It is clean and quiet. However there were BED BUGS in room 341. The other people that stay there are a reasonable class. Without the bugs it is a good place to stay.'
Set-NumCell $ws2 4 13 3
Set-NumCell $ws2 4 16 4
Set-NumCell $ws2 4 17 2
Set-NumCell $ws2 4 18 4
Set-NumCell $ws2 4 19 2
Set-NumCell $ws2 4 21 3
Set-NumCell $ws2 4 22 0
Set-TextCell $ws2 4 25 'I lived here 6 months. It is a great location---the management is nice. It is clean and quiet. However there were BED BUGS in room 341. The other people that stay there are a reasonable class. Without the bugs it is a good place to stay.'

Write-Host "edit applied"
